$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated CRLB formula now divides by an additional factor of 2^2 (4),
# effectively halving every CRLB [Hz^2] value in column F (rows 2-49).
# The Monte Carlo simulation was re-run, which also updated the
# "Mean estimated" (C), "Mean error" (D) and "Variance" (E) columns
# for the rows where the estimator had not yet converged to the true value.

$ws.Range("C2").Value = 99787.109375
$ws.Range("D2").Value = 212.890625
$ws.Range("E2").Value = 142121.338868165
$ws.Range("F2").Value = 11257.47700647604
$ws.Range("F3").Value = 1125.747700647604
$ws.Range("F4").Value = 112.5747700647604
$ws.Range("F5").Value = 11.25747700647604
$ws.Range("F6").Value = 1.125747700647604
$ws.Range("F7").Value = 0.1125747700647604
$ws.Range("F8").Value = 0.01125747700647604
$ws.Range("F9").Value = 0.001125747700647604
$ws.Range("C10").Value = 99991.69921875
$ws.Range("D10").Value = 8.30078125
$ws.Range("E10").Value = 19309.99564933586
$ws.Range("F10").Value = 11257.47700647604
$ws.Range("C11").Value = 100040.771484375
$ws.Range("D11").Value = -40.771484375
$ws.Range("E11").Value = 10662.66834079563
$ws.Range("F11").Value = 1125.747700647604
$ws.Range("C12").Value = 100094.970703125
$ws.Range("D12").Value = -94.970703125
$ws.Range("E12").Value = 649.0880185300047
$ws.Range("F12").Value = 112.5747700647604
$ws.Range("F13").Value = 11.25747700647604
$ws.Range("F14").Value = 1.125747700647604
$ws.Range("F15").Value = 0.1125747700647604
$ws.Range("F16").Value = 0.01125747700647604
$ws.Range("F17").Value = 0.001125747700647604
$ws.Range("C18").Value = 99999.51171875
$ws.Range("D18").Value = 0.48828125
$ws.Range("E18").Value = 11449.34225607443
$ws.Range("F18").Value = 11257.47700647604
$ws.Range("C19").Value = 100000.0610351562
$ws.Range("D19").Value = -0.06103515625
$ws.Range("E19").Value = 1552.014111040591
$ws.Range("F19").Value = 1125.747700647604
$ws.Range("C20").Value = 99993.5302734375
$ws.Range("D20").Value = 6.4697265625
$ws.Range("E20").Value = 774.0101656756243
$ws.Range("F20").Value = 112.5747700647604
$ws.Range("C21").Value = 99977.72216796875
$ws.Range("D21").Value = 22.27783203125
$ws.Range("E21").Value = 125.9476274580092
$ws.Range("F21").Value = 11.25747700647604
$ws.Range("F22").Value = 1.125747700647604
$ws.Range("F23").Value = 0.1125747700647604
$ws.Range("F24").Value = 0.01125747700647604
$ws.Range("F25").Value = 0.001125747700647604
$ws.Range("C26").Value = 100003.4027099609
$ws.Range("D26").Value = -3.4027099609375
$ws.Range("E26").Value = 11330.08239169916
$ws.Range("F26").Value = 11257.47700647604
$ws.Range("C27").Value = 100000.5187988281
$ws.Range("D27").Value = -0.518798828125
$ws.Range("E27").Value = 1192.830309077903
$ws.Range("F27").Value = 1125.747700647604
$ws.Range("C28").Value = 99999.42016601562
$ws.Range("D28").Value = 0.579833984375
$ws.Range("E28").Value = 133.3487985489724
$ws.Range("F28").Value = 112.5747700647604
$ws.Range("C29").Value = 100001.1444091797
$ws.Range("D29").Value = -1.1444091796875
$ws.Range("E29").Value = 51.12835080237002
$ws.Range("F29").Value = 11.25747700647604
$ws.Range("C30").Value = 100004.7607421875
$ws.Range("D30").Value = -4.7607421875
$ws.Range("E30").Value = 18.70476089798294
$ws.Range("F30").Value = 1.125747700647604
$ws.Range("F31").Value = 0.1125747700647604
$ws.Range("F32").Value = 0.01125747700647604
$ws.Range("F33").Value = 0.001125747700647604
$ws.Range("C34").Value = 100006.7443847656
$ws.Range("D34").Value = -6.744384765625
$ws.Range("E34").Value = 11621.28151022457
$ws.Range("F34").Value = 11257.47700647604
$ws.Range("C35").Value = 100001.8005371094
$ws.Range("D35").Value = -1.800537109375
$ws.Range("E35").Value = 1079.254169334043
$ws.Range("F35").Value = 1125.747700647604
$ws.Range("C36").Value = 100000.3280639648
$ws.Range("D36").Value = -0.32806396484375
$ws.Range("E36").Value = 115.8006745646964
$ws.Range("F36").Value = 112.5747700647604
$ws.Range("C37").Value = 100000.1640319824
$ws.Range("D37").Value = -0.164031982421875
$ws.Range("E37").Value = 12.53811369843043
$ws.Range("F37").Value = 11.25747700647604
$ws.Range("C38").Value = 99999.80545043945
$ws.Range("D38").Value = 0.194549560546875
$ws.Range("E38").Value = 3.338623041467922
$ws.Range("F38").Value = 1.125747700647604
$ws.Range("C39").Value = 99998.96240234375
$ws.Range("D39").Value = 1.03759765625
$ws.Range("E39").Value = 1.625852422551946
$ws.Range("F39").Value = 0.1125747700647604
$ws.Range("F40").Value = 0.01125747700647604
$ws.Range("F41").Value = 0.001125747700647604
$ws.Range("C42").Value = 100000.1211166382
$ws.Range("D42").Value = -0.1211166381835938
$ws.Range("E42").Value = 11253.58350410981
$ws.Range("F42").Value = 11257.47700647604
$ws.Range("C43").Value = 100000.5912780762
$ws.Range("D43").Value = -0.591278076171875
$ws.Range("E43").Value = 1169.919360054282
$ws.Range("F43").Value = 1125.747700647604
$ws.Range("C44").Value = 99999.47071075439
$ws.Range("D44").Value = 0.5292892456054688
$ws.Range("E44").Value = 109.8631135125932
$ws.Range("F44").Value = 112.5747700647604
$ws.Range("C45").Value = 99999.84645843506
$ws.Range("D45").Value = 0.1535415649414062
$ws.Range("E45").Value = 11.489566451848
$ws.Range("F45").Value = 11.25747700647604
$ws.Range("C46").Value = 100000.0305175781
$ws.Range("D46").Value = -0.030517578125
$ws.Range("E46").Value = 1.165842930147717
$ws.Range("F46").Value = 1.125747700647604
$ws.Range("C47").Value = 100000.036239624
$ws.Range("D47").Value = -0.0362396240234375
$ws.Range("E47").Value = 0.2230091933608503
$ws.Range("F47").Value = 0.1125747700647604
$ws.Range("C48").Value = 100000.2145767212
$ws.Range("D48").Value = -0.2145767211914062
$ws.Range("E48").Value = 0.1314397373057723
$ws.Range("F48").Value = 0.01125747700647604
$ws.Range("C49").Value = 100000.3786087036
$ws.Range("D49").Value = -0.3786087036132812
$ws.Range("E49").Value = 0.002723021674677506
$ws.Range("F49").Value = 0.001125747700647604
